$d = $word.ActiveDocument

# 1) Fix the date in the first paragraph: 26.07.24 -> 25.07.24
$d.Content.Find.Execute("26.07.24", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "25.07.24", 2) | Out-Null

# 2) Insert a new Heading1 paragraph right after paragraph 1 (the date line)
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter() | Out-Null
$headingPara = $d.Paragraphs(2)
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "AI models collapse when trained on recursively generated data"

# 3) Replace the old "Questionable practices..." paragraph (now paragraph 3)
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "מאמר די חמוד שחוקר מה קורה שמאמנים מודלי AI על הדאטה הנוצר על ידי מודלי AI. בשתי מילים - לא הכל ורוד שם ויש כמה סיבות למה הדברים עלולים להשתבש:"

# 4) Replace the long intro paragraph (now paragraph 4). This paragraph's run
#    originally ended with a trailing space (xml:space="preserve"); reusing the
#    same run would keep that stale flag, so insert a brand-new paragraph with
#    the new text and delete the old tainted one instead of overwriting in place.
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertParagraphAfter() | Out-Null
$newP4 = $d.Paragraphs(5)
$newP4.Range.Text = "דאטה דריפט (איך זה בעברית?) קיצוני: אימון מודלים על דאטה שנוצרה על ידי מודלים אחרים גורם להתרחקות של התפלגות הדאטה הנוצר על ידי המודל החדש מהדאטה האמיתי (כלומר אגרגציה של מרחק בין ההתפלגויות שלהן).."
$oldP4 = $d.Paragraphs(4)
$oldP4.Range.Delete() | Out-Null

# 5) Replace the URL paragraph (now paragraph 5) with the next bullet point of text
$p5 = $d.Paragraphs(5)
$p5.Range.Text = "הבעיות מחמירות בזנבות התפלגות הדאטה (תחומים או שפות עם מעט דאטה למשל): ההידרדרות משפיעה בעיקר על זנבות התפלגות הדאטה, שם דאטה נדיר הופך להיות עוד פחות מיוצג"

# 6) Append three more paragraphs after paragraph 5: two more bullet points and the new URL
$r5 = $d.Paragraphs(5).Range
$r5.Collapse(0)
$r5.InsertParagraphAfter() | Out-Null
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "עוד יותר שגיאות: שגיאות בדאטה שנוצרו על ידי מודלים מצטברות לאורך דורות, מה שמוביל לירידה משמעותית בביצועים."

$r6 = $d.Paragraphs(6).Range
$r6.Collapse(0)
$r6.InsertParagraphAfter() | Out-Null
$p7 = $d.Paragraphs(7)
$p7.Range.Text = "קריסת השונות: דאטה שנוצר על ידי מודלים חסרים את המגוון והעושר של הדאטה מהעולם האמיתי, מה שמוביל ליותר הומוגניזציית יתר (פחות גיוון)."

$r7 = $d.Paragraphs(7).Range
$r7.Collapse(0)
$r7.InsertParagraphAfter() | Out-Null
$p8 = $d.Paragraphs(8)
$p8.Range.Text = "https://www.nature.com/articles/s41586-024-07566-y"

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
